$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '2026-02-21 23:18:17'
$ws.Range("O2").Value = '3.7 °C'
$ws.Range("E3").Value = '2026-02-21 23:18:20'
$ws.Range("O3").Value = '2.1 °C'
$ws.Range("E4").Value = '2026-02-21 23:18:22'
$ws.Range("K4").Value = '14.5 MJ/m2'
$ws.Range("O4").Value = '9.0 °C'
$ws.Range("E5").Value = '2026-02-21 23:18:25'
$ws.Range("O5").Value = '3.9 °C'
$ws.Range("E6").Value = '2026-02-21 23:18:27'
$ws.Range("E7").Value = '2026-02-21 23:18:29'
$ws.Range("O7").Value = '13.3 °C'
$ws.Range("E8").Value = '2026-02-21 23:18:32'
$ws.Range("H8").NumberFormat = "@"
$ws.Range("H8").Value = '61%'
$ws.Range("E9").Value = '2026-02-21 23:18:34'
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = '60%'
$ws.Range("O9").Value = '12.7 °C'
$ws.Range("E10").Value = '2026-02-21 23:18:37'
$ws.Range("H10").NumberFormat = "@"
$ws.Range("H10").Value = '80%'
$ws.Range("K10").Value = '14.8 MJ/m2'
$ws.Range("O10").Value = '8.2 °C'
$ws.Range("E11").Value = '2026-02-21 23:18:39'
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H11").Value = '55%'
$ws.Range("O11").Value = '8.4 °C'
$ws.Range("E12").Value = '2026-02-21 23:18:42'
$ws.Range("H12").NumberFormat = "@"
$ws.Range("H12").Value = '66%'
$ws.Range("N12").Value = '6.3 °C 22:53 TU'
$ws.Range("O12").Value = '12.2 °C'
$ws.Range("E13").Value = '2026-02-21 23:18:44'
$ws.Range("J13").Value = '1032.1 hPa'
$ws.Range("O13").Value = '5.1 °C'
$ws.Range("E14").Value = '2026-02-21 23:18:46'
$ws.Range("O14").Value = '11.0 °C'
$ws.Range("E15").Value = '2026-02-21 23:18:49'
$ws.Range("H15").NumberFormat = "@"
$ws.Range("H15").Value = '59%'
$ws.Range("O15").Value = '12.5 °C'
$ws.Range("E16").Value = '2026-02-21 23:18:51'
$ws.Range("O16").Value = '2.4 °C'
$ws.Range("E17").Value = '2026-02-21 23:18:53'
$ws.Range("O17").Value = '8.6 °C'
$ws.Range("E18").Value = '2026-02-21 23:18:56'
$ws.Range("O18").Value = '8.3 °C'
$ws.Range("E19").Value = '2026-02-21 23:18:58'
$ws.Range("E20").Value = '2026-02-21 23:19:01'
$ws.Range("E21").Value = '2026-02-21 23:19:03'
$ws.Range("J21").Value = '1030.9 hPa'
$ws.Range("O21").Value = '7.5 °C'
$ws.Range("E22").Value = '2026-02-21 23:19:05'
$ws.Range("E23").Value = '2026-02-21 23:19:08'
$ws.Range("E24").Value = '2026-02-21 23:19:10'
$ws.Range("J24").Value = '1031.7 hPa'
$ws.Range("O24").Value = '6.1 °C'
$ws.Range("E25").Value = '2026-02-21 23:19:13'
$ws.Range("E26").Value = '2026-02-21 23:19:15'
$ws.Range("E27").Value = '2026-02-21 23:19:17'
$ws.Range("H27").NumberFormat = "@"
$ws.Range("H27").Value = '33%'
$ws.Range("K27").Value = '16.4 MJ/m2'
$ws.Range("E28").Value = '2026-02-21 23:19:20'
$ws.Range("E29").Value = '2026-02-21 23:19:22'
$ws.Range("O29").Value = '11.1 °C'
$ws.Range("E30").Value = '2026-02-21 23:19:25'
$ws.Range("O30").Value = '11.2 °C'
$ws.Range("E31").Value = '2026-02-21 23:19:27'
$ws.Range("E32").Value = '2026-02-21 23:19:30'
$ws.Range("H32").NumberFormat = "@"
$ws.Range("H32").Value = '82%'
$ws.Range("O32").Value = '4.5 °C'
$ws.Range("E33").Value = '2026-02-21 23:19:32'
$ws.Range("J33").Value = '1030.6 hPa'
$ws.Range("E34").Value = '2026-02-21 23:19:34'
$ws.Range("E35").Value = '2026-02-21 23:19:37'
$ws.Range("H35").NumberFormat = "@"
$ws.Range("H35").Value = '53%'
$ws.Range("J35").Value = '1031.2 hPa'
$ws.Range("O35").Value = '7.4 °C'
$ws.Range("E36").Value = '2026-02-21 23:19:39'
$ws.Range("H36").NumberFormat = "@"
$ws.Range("H36").Value = '61%'
$ws.Range("O36").Value = '12.8 °C'
$ws.Range("E37").Value = '2026-02-21 23:19:42'
$ws.Range("J37").Value = '1031.8 hPa'
$ws.Range("O37").Value = '5.4 °C'
$ws.Range("E38").Value = '2026-02-21 23:19:44'
$ws.Range("O38").Value = '9.5 °C'
$ws.Range("E39").Value = '2026-02-21 23:19:46'
$ws.Range("E40").Value = '2026-02-21 23:19:49'
$ws.Range("J40").Value = '1030.9 hPa'
$ws.Range("O40").Value = '8.2 °C'
$ws.Range("E41").Value = '2026-02-21 23:19:51'
$ws.Range("H41").NumberFormat = "@"
$ws.Range("H41").Value = '71%'
$ws.Range("O41").Value = '11.0 °C'
$ws.Range("E42").Value = '2026-02-21 23:19:54'
$ws.Range("H42").NumberFormat = "@"
$ws.Range("H42").Value = '77%'
$ws.Range("O42").Value = '10.5 °C'
$ws.Range("E43").Value = '2026-02-21 23:19:56'
$ws.Range("E44").Value = '2026-02-21 23:19:59'
$ws.Range("E45").Value = '2026-02-21 23:20:01'
$ws.Range("E46").Value = '2026-02-21 23:20:03'
$ws.Range("H46").NumberFormat = "@"
$ws.Range("H46").Value = '71%'
$ws.Range("O46").Value = '9.4 °C'
